# Insert a new weekly price record as row 66 in the "Hortaliza, Feria Lagunitas
# de Puerto Montt - Zapallo" sheet. This pushes the existing rows 66-149 down to
# 67-150 (the sheet's used range grows from A1:R149 to A1:R150), and fills the
# newly created row 66 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above current row 66 (shifts rows 66..149 -> 67..150).
$ws.Rows.Item(66).Insert()

# Populate the new row 66 with the new record.
$ws.Cells.Item(66, 1).Value  = 4
$ws.Cells.Item(66, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(66, 3).Value  = "Los Lagos"
$ws.Cells.Item(66, 4).Value  = 44413
$ws.Cells.Item(66, 5).Value  = 10
$ws.Cells.Item(66, 6).Value  = 100112045
$ws.Cells.Item(66, 7).Value  = "Zapallo"
$ws.Cells.Item(66, 8).Value  = "Paine"
$ws.Cells.Item(66, 9).Value  = "1a (guarda)"
$ws.Cells.Item(66, 10).Value = 500
$ws.Cells.Item(66, 11).Value = 450
$ws.Cells.Item(66, 12).Value = 450
$ws.Cells.Item(66, 13).Value = 450
$ws.Cells.Item(66, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(66, 15).Value = "Región Metropolitana"
$ws.Cells.Item(66, 16).Value = 450
$ws.Cells.Item(66, 17).Value = 1
$ws.Cells.Item(66, 18).Value = "Hortaliza"
